$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.2315112540192926
$ws.Range("C2").Value = 0.4694533762057878
$ws.Range("J2").Value = 0.01607717041800643
$ws.Range("P2").Value = 0.1414790996784566
$ws.Range("S2").Value = 0.1414790996784566
$ws.Range("B3").Value = 0.006535947712418301
$ws.Range("C3").Value = 0.0392156862745098
$ws.Range("J3").Value = 0.0718954248366013
$ws.Range("P3").Value = 0.6666666666666666
$ws.Range("S3").Value = 0.2156862745098039
$ws.Range("J4").Value = 0.1304347826086956
$ws.Range("P4").Value = 0.6521739130434783
$ws.Range("S4").Value = 0.2173913043478261
$ws.Range("B6").Value = 0.07692307692307693
$ws.Range("D6").Value = 0.01923076923076923
$ws.Range("F6").Value = 0.06538461538461539
$ws.Range("J6").Value = 0.2115384615384615
$ws.Range("O6").Value = 0.03076923076923077
$ws.Range("Q6").Value = 0.1153846153846154
$ws.Range("R6").Value = 0.04230769230769231
$ws.Range("S6").Value = 0.4384615384615385
$ws.Range("B7").Value = 0.1428571428571428
$ws.Range("D7").Value = 0.01020408163265306
$ws.Range("E7").Value = 0.00510204081632653
$ws.Range("F7").Value = 0.08673469387755102
$ws.Range("J7").Value = 0.09183673469387756
$ws.Range("O7").Value = 0.04081632653061224
$ws.Range("Q7").Value = 0.1173469387755102
$ws.Range("R7").Value = 0.1122448979591837
$ws.Range("S7").Value = 0.3928571428571428
$ws.Range("B8").Value = 0.09640831758034027
$ws.Range("D8").Value = 0.02079395085066163
$ws.Range("E8").Value = 0.001890359168241966
$ws.Range("F8").Value = 0.08506616257088846
$ws.Range("J8").Value = 0.1266540642722117
$ws.Range("O8").Value = 0.00945179584120983
$ws.Range("Q8").Value = 0.1474480151228733
$ws.Range("R8").Value = 0.0888468809073724
$ws.Range("S8").Value = 0.4234404536862004
$ws.Range("B9").Value = 0.118942731277533
$ws.Range("D9").Value = 0.013215859030837
$ws.Range("F9").Value = 0.0881057268722467
$ws.Range("J9").Value = 0.1101321585903084
$ws.Range("O9").Value = 0.01762114537444934
$ws.Range("Q9").Value = 0.1718061674008811
$ws.Range("R9").Value = 0.07488986784140969
$ws.Range("S9").Value = 0.4052863436123348
$ws.Range("B10").Value = 0.08592592592592592
$ws.Range("D10").Value = 0.02074074074074074
$ws.Range("E10").Value = 0.001481481481481481
$ws.Range("F10").Value = 0.07333333333333333
$ws.Range("J10").Value = 0.12
$ws.Range("O10").Value = 0.01703703703703704
$ws.Range("Q10").Value = 0.2014814814814815
$ws.Range("R10").Value = 0.08666666666666667
$ws.Range("S10").Value = 0.3933333333333333
$ws.Range("G11").Value = 0.1656626506024096
$ws.Range("J11").Value = 0.09036144578313253
$ws.Range("K11").Value = 0.2078313253012048
$ws.Range("L11").Value = 0.5150602409638554
$ws.Range("S11").Value = 0.02108433734939759
$ws.Range("G12").Value = 0.6820809248554913
$ws.Range("J12").Value = 0.2427745664739884
$ws.Range("K12").Value = 0.005780346820809248
$ws.Range("L12").Value = 0.01734104046242774
$ws.Range("S12").Value = 0.05202312138728324
$ws.Range("G13").Value = 0.717948717948718
$ws.Range("J13").Value = 0.282051282051282
$ws.Range("F15").Value = 0.01158301158301158
$ws.Range("H15").Value = 0.1698841698841699
$ws.Range("I15").Value = 0.08494208494208494
$ws.Range("J15").Value = 0.4131274131274131
$ws.Range("K15").Value = 0.07335907335907337
$ws.Range("M15").Value = 0.007722007722007722
$ws.Range("O15").Value = 0.05791505791505792
$ws.Range("S15").Value = 0.1814671814671815
$ws.Range("F16").Value = 0.02352941176470588
$ws.Range("H16").Value = 0.1941176470588235
$ws.Range("I16").Value = 0.07058823529411765
$ws.Range("J16").Value = 0.3764705882352941
$ws.Range("K16").Value = 0.1352941176470588
$ws.Range("M16").Value = 0.02352941176470588
$ws.Range("O16").Value = 0.05882352941176471
$ws.Range("S16").Value = 0.1176470588235294
$ws.Range("F17").Value = 0.01601830663615561
$ws.Range("H17").Value = 0.2036613272311213
$ws.Range("I17").Value = 0.1052631578947368
$ws.Range("J17").Value = 0.414187643020595
$ws.Range("K17").Value = 0.08466819221967964
$ws.Range("M17").Value = 0.02288329519450801
$ws.Range("N17").Value = 0.002288329519450801
$ws.Range("O17").Value = 0.05720823798627003
$ws.Range("S17").Value = 0.09382151029748284
$ws.Range("F18").Value = 0.02336448598130841
$ws.Range("H18").Value = 0.2383177570093458
$ws.Range("I18").Value = 0.1074766355140187
$ws.Range("J18").Value = 0.3598130841121495
$ws.Range("K18").Value = 0.1214953271028037
$ws.Range("M18").Value = 0.02803738317757009
$ws.Range("N18").Value = 0.004672897196261682
$ws.Range("O18").Value = 0.07943925233644859
$ws.Range("S18").Value = 0.03738317757009346
$ws.Range("F19").Value = 0.01711840228245364
$ws.Range("H19").Value = 0.2232524964336662
$ws.Range("I19").Value = 0.0898716119828816
$ws.Range("J19").Value = 0.3594864479315264
$ws.Range("K19").Value = 0.1105563480741797
$ws.Range("M19").Value = 0.01283880171184023
$ws.Range("O19").Value = 0.07703281027104136
$ws.Range("S19").Value = 0.1098430813124108
